# Update the "5.3" sheet: two data values were corrected (they were
# mistakenly entered without a decimal point) and that sheet became the
# active/selected tab (with a new active-cell selection), replacing "5.2"
# as the previously selected tab.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("5.3")
$ws.Activate()

$ws.Range("C7").Value = 7.24
$ws.Range("C8").Value = 7.88

$ws.Range("C9").Select()
